$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '60.677.60'
$c.ClearFormats()
$ws.Range("E2").Value = '  +2.76%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.703.99'
$c.ClearFormats()
$ws.Range("E3").Value = '  +2.57%  '

$ws.Range("E4").Value = '  +0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '526.58'
$c.ClearFormats()
$ws.Range("E5").Value = '  +1.60%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '145.05'
$c.ClearFormats()
$ws.Range("E6").Value = '  -0.48%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +1.80%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.732.80'
$c.ClearFormats()
$ws.Range("E9").Value = '  +2.66%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '6.71'
$c.ClearFormats()
$ws.Range("E10").Value = '  +6.94%  '

$ws.Range("E11").Value = '  +0.86%  '

$ws.Range("E12").Value = '  +0.86%  '

$ws.Range("E13").Value = '  +2.93%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '3.182.55'
$c.ClearFormats()
$ws.Range("E14").Value = '  +2.51%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '60.645.49'
$c.ClearFormats()
$ws.Range("E15").Value = '  +2.77%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.860.32'
$c.ClearFormats()
$ws.Range("E16").Value = '  +7.63%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '21.35'
$c.ClearFormats()
$ws.Range("E17").Value = '  +1.51%  '

$ws.Range("E18").Value = '  +0.49%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '348.68'
$c.ClearFormats()
$ws.Range("E19").Value = '  -0.48%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.52'
$c.ClearFormats()
$ws.Range("E20").Value = '  -0.43%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '10.60'
$c.ClearFormats()
$ws.Range("E21").Value = '  +2.11%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.40'
$c.ClearFormats()
$ws.Range("E22").Value = '  +3.13%  '

$ws.Range("E23").Value = '  -0.33%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '63.68'
$c.ClearFormats()
$ws.Range("E24").Value = '  +2.95%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.420'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("E26").Value = '  +4.72%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.ClearFormats()
$ws.Range("E27").Value = '  -0.29%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.0₃0819'
$c.ClearFormats()
$ws.Range("E28").Value = '  +1.17%  '

$ws.Range("E29").Value = '  +2.06%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.77'
$c.ClearFormats()
$ws.Range("E30").Value = '  +7.68%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range("E31").Value = '  +0.04%  '

$ws.Range("E32").Value = '  +1.67%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '19.16'
$c.ClearFormats()
$ws.Range("E33").Value = '  +0.72%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '150.62'
$c.ClearFormats()
$ws.Range("E34").Value = '  +0.26%  '

$ws.Range("E35").Value = '  +5.28%  '

$ws.Range("E36").Value = '  +7.65%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.947'
$c.ClearFormats()
$ws.Range("E37").Value = '  -1.90%  '

$ws.Range("E38").Value = '  +3.76%  '

$ws.Range("E39").Value = '  +7.51%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '37.01'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.61%  '

$ws.Range("E41").Value = '  -1.23%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '285.10'
$c.ClearFormats()
$ws.Range("E42").Value = '  +2.50%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '20.16'
$c.ClearFormats()
$ws.Range("E43").Value = '  +2.72%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.614'
$c.ClearFormats()
$ws.Range("E44").Value = '  +0.61%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0992'
$c.ClearFormats()
$ws.Range("E45").Value = '  +0.64%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.148.59'
$c.ClearFormats()
$ws.Range("E46").Value = '  +7.60%  '

$ws.Range("E47").Value = '  +0.10%  '

$ws.Range("E48").Value = '  +2.64%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0235'
$c.ClearFormats()
$ws.Range("E49").Value = '  +2.09%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '4.81'
$c.ClearFormats()
$ws.Range("E50").Value = '  +1.85%  '

$ws.Range("E51").Value = '  +1.65%  '
